$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2"); $c.Value = "'26.331.26"; $c.Style = "Normal"
$ws.Range("E2").Value = "  -4.18%  "
$c = $ws.Range("D3"); $c.Value = "'1.761.75"; $c.Style = "Normal"
$ws.Range("E3").Value = "  -3.29%  "
$ws.Range("E4").Value = "  +0.02%  "
$c = $ws.Range("D6"); $c.Value = "'304.32"; $c.Style = "Normal"
$ws.Range("E6").Value = "  -2.37%  "
$c = $ws.Range("D7"); $c.Value = "'0.4271"; $c.Style = "Normal"
$ws.Range("E7").Value = "  +0.83%  "
$c = $ws.Range("D8"); $c.Value = "'0.3607"; $c.Style = "Normal"
$ws.Range("E8").Value = "  -0.48%  "
$c = $ws.Range("D9"); $c.Value = "'0.07052"; $c.Style = "Normal"
$ws.Range("E9").Value = "  -1.73%  "
$c = $ws.Range("D10"); $c.Value = "'0.8311"; $c.Style = "Normal"
$ws.Range("E10").Value = "  -3.04%  "
$c = $ws.Range("D11"); $c.Value = "'20.13"; $c.Style = "Normal"
$ws.Range("E11").Value = "  -2.10%  "
$c = $ws.Range("D12"); $c.Value = "'1.793.95"; $c.Style = "Normal"
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("E13").Value = "  -3.09%  "
$c = $ws.Range("D14"); $c.Value = "'6.398"; $c.Style = "Normal"
$ws.Range("E14").Value = "  -0.99%  "
$c = $ws.Range("D15"); $c.Value = "'0.06813"; $c.Style = "Normal"
$ws.Range("E15").Value = "  -1.40%  "
$c = $ws.Range("D16"); $c.Value = "'1.006"; $c.Style = "Normal"
$ws.Range("E16").Value = "  +0.12%  "
$c = $ws.Range("D17"); $c.Value = "'79.15"; $c.Style = "Normal"
$ws.Range("E17").Value = "  -1.18%  "
$c = $ws.Range("D18"); $c.Value = "'0.000008624"; $c.Style = "Normal"
$ws.Range("E18").Value = "  -2.66%  "
$c = $ws.Range("D19"); $c.Value = "'1.002"; $c.Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$c = $ws.Range("D20"); $c.Value = "'14.95"; $c.Style = "Normal"
$ws.Range("E20").Value = "  -2.29%  "
$c = $ws.Range("D21"); $c.Value = "'26.351.47"; $c.Style = "Normal"
$ws.Range("E21").Value = "  -4.84%  "
$c = $ws.Range("D22"); $c.Value = "'4.993"; $c.Style = "Normal"
$ws.Range("E22").Value = "  -2.32%  "
$c = $ws.Range("D23"); $c.Value = "'11.10"; $c.Style = "Normal"
$ws.Range("E23").Value = "  +1.95%  "
$c = $ws.Range("D24"); $c.Value = "'1.990.37"; $c.Style = "Normal"
$ws.Range("E24").Value = "  -1.74%  "
$c = $ws.Range("D25"); $c.Value = "'1.892"; $c.Style = "Normal"
$ws.Range("E25").Value = "  -4.77%  "
$c = $ws.Range("D26"); $c.Value = "'152.47"; $c.Style = "Normal"
$ws.Range("E26").Value = "  -1.63%  "
$c = $ws.Range("D27"); $c.Value = "'18.09"; $c.Style = "Normal"
$ws.Range("E27").Value = "  -3.16%  "
$c = $ws.Range("D28"); $c.Value = "'114.59"; $c.Style = "Normal"
$ws.Range("E28").Value = "  +0.64%  "
$c = $ws.Range("D29"); $c.Value = "'5.018"; $c.Style = "Normal"
$ws.Range("E29").Value = "  -2.25%  "
$c = $ws.Range("D30"); $c.Value = "'1.680"; $c.Style = "Normal"
$c = $ws.Range("D31"); $c.Value = "'0.08890"; $c.Style = "Normal"
$ws.Range("E31").Value = "  +0.62%  "
$c = $ws.Range("D32"); $c.Value = "'0.7217"; $c.Style = "Normal"
$ws.Range("E32").Value = "  -2.70%  "
$c = $ws.Range("D33"); $c.Value = "'4.309"; $c.Style = "Normal"
$ws.Range("E33").Value = "  -4.63%  "
$c = $ws.Range("D34"); $c.Value = "'1.102"; $c.Style = "Normal"
$ws.Range("E34").Value = "  -1.50%  "
$c = $ws.Range("D35"); $c.Value = "'2.752"; $c.Style = "Normal"
$ws.Range("E35").Value = "  -7.73%  "
$c = $ws.Range("D36"); $c.Value = "'1.002"; $c.Style = "Normal"
$ws.Range("E36").Value = "  +0.00%  "
$c = $ws.Range("D37"); $c.Value = "'1.068"; $c.Style = "Normal"
$ws.Range("E37").Value = "  -1.66%  "
$c = $ws.Range("D38"); $c.Value = "'0.05088"; $c.Style = "Normal"
$ws.Range("E38").Value = "  -3.42%  "
$ws.Range("E39").Value = "  -1.98%  "
$c = $ws.Range("D40"); $c.Value = "'0.4890"; $c.Style = "Normal"
$ws.Range("E40").Value = "  -2.93%  "
$c = $ws.Range("D41"); $c.Value = "'0.1596"; $c.Style = "Normal"
$ws.Range("E41").Value = "  -2.61%  "
$c = $ws.Range("D42"); $c.Value = "'6.204"; $c.Style = "Normal"
$ws.Range("E42").Value = "  -3.83%  "
$c = $ws.Range("D43"); $c.Value = "'2.498"; $c.Style = "Normal"
$ws.Range("E43").Value = "  -10.29%  "
$c = $ws.Range("D44"); $c.Value = "'7.979"; $c.Style = "Normal"
$ws.Range("E44").Value = "  -3.64%  "
$c = $ws.Range("D45"); $c.Value = "'104.72"; $c.Style = "Normal"
$ws.Range("E45").Value = "  -0.72%  "
$c = $ws.Range("D46"); $c.Value = "'1.002"; $c.Style = "Normal"
$ws.Range("E46").Value = "  +0.11%  "
$c = $ws.Range("D47"); $c.Value = "'10.04"; $c.Style = "Normal"
$ws.Range("E47").Value = "  -2.97%  "
$c = $ws.Range("D49"); $c.Value = "'0.4467"; $c.Style = "Normal"
$ws.Range("E49").Value = "  -3.99%  "
$c = $ws.Range("D50"); $c.Value = "'1.568"; $c.Style = "Normal"
$ws.Range("E50").Value = "  -2.47%  "
$c = $ws.Range("D51"); $c.Value = "'1.714"; $c.Style = "Normal"
$ws.Range("E51").Value = "  -0.18%  "
